$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new image URL values on row 2 (columns G, H, I) ---
$ws.Range("G2").Value = "https://storagecencosud.blob.core.windows.net/nathaly/producto.png"
$ws.Range("H2").Value = "https://storagecencosud.blob.core.windows.net/nathaly/material.png"
$ws.Range("I2").Value = "https://storagecencosud.blob.core.windows.net/nathaly/dimensiones.png"

# --- Bold header row A1:I1 ---
$ws.Range("A1:I1").Font.Bold = $true

# --- Header fill highlight on G1:I1 ---
$ws.Range("G1:I1").Interior.ThemeColor = 9
$ws.Range("G1:I1").Interior.TintAndShade = 0.79998168889431442

# --- Row height / default row height ---
$ws.Range("A1:I5").RowHeight = 15.5

# --- Column widths ---
$ws.Columns.Item(7).ColumnWidth = 31.4140625
$ws.Columns.Item(8).ColumnWidth = 24.75

# --- Selection ---
$ws.Range("E8").Select()
